$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 15877.889
$ws.Range("J13").Value = 15877.889
$ws.Range("L13").Value = 15877.889
$ws.Range("N13").Value = -16215.889

$ws.Range("H17").Value = 3244.25
$ws.Range("J17").Value = 3292.851
$ws.Range("L17").Value = 9878.553
$ws.Range("N17").Value = -10214.553

$ws.Range("H18").Value = 737.75
$ws.Range("I18").Value = 737.75
$ws.Range("K18").Value = 737.75
$ws.Range("M18").Value = -453.75

$ws.Range("H28").Value = 372
$ws.Range("I28").Value = 356.64285
$ws.Range("J28").Value = 443.66666
$ws.Range("K28").Value = 356.64285
$ws.Range("L28").Value = 443.66666
$ws.Range("M28").Value = 128.35715
$ws.Range("N28").Value = -1413.66666

$ws.Range("H32").Value = 14493552
$ws.Range("I32").Value = 41666948
$ws.Range("J32").Value = 1073.8667
$ws.Range("K32").Value = 41666948
$ws.Range("L32").Value = 1073.8667
$ws.Range("M32").Value = -41666622
$ws.Range("N32").Value = -1725.8667

$ws.Range("H64").Value = 26437.904
$ws.Range("I64").Value = 168916.67
$ws.Range("J64").Value = 2691.4443
$ws.Range("K64").Value = 168916.67
$ws.Range("L64").Value = 2691.4443
$ws.Range("M64").Value = -168668.67
$ws.Range("N64").Value = -3187.4443

$ws.Range("H67").Value = 26437.904
$ws.Range("I67").Value = 168916.67
$ws.Range("J67").Value = 2691.4443
$ws.Range("K67").Value = 168916.67
$ws.Range("L67").Value = 2691.4443
$ws.Range("M67").Value = -168058.67
$ws.Range("N67").Value = -4407.4443

$ws.Range("H114").Value = 41198
$ws.Range("J114").Value = 41198
$ws.Range("L114").Value = 41198
$ws.Range("N114").Value = -49876

$ws.Range("H127").Value = 1440.5714
$ws.Range("I127").Value = 559.6
$ws.Range("J127").Value = 1930
$ws.Range("K127").Value = 1678.8
$ws.Range("L127").Value = 5790
$ws.Range("M127").Value = 3281.2
$ws.Range("N127").Value = -15710

$ws.Range("H132").Value = 22097.723
$ws.Range("I132").Value = 3560.889
$ws.Range("J132").Value = 82763.73
$ws.Range("K132").Value = 10682.667
$ws.Range("L132").Value = 248291.19
$ws.Range("M132").Value = -8152.667000000001
$ws.Range("N132").Value = -253351.19

$ws.Range("H135").Value = 15626365
$ws.Range("I135").Value = 1392.7778
$ws.Range("J135").Value = 35715616
$ws.Range("K135").Value = 12535.0002
$ws.Range("L135").Value = 321440544
$ws.Range("M135").Value = -10000.0002
$ws.Range("N135").Value = -321445614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 8000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 8000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 8000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -8598

$ws.Range("H32").Value = 42306.54
$ws.Range("I32").Value = 43312.6
$ws.Range("J32").Value = 33503.5
$ws.Range("K32").Value = 43312.6
$ws.Range("L32").Value = 33503.5
$ws.Range("M32").Value = -43025.6
$ws.Range("N32").Value = -34077.5

$ws.Range("H75").Value = 25586.5
$ws.Range("J75").Value = 25586.5
$ws.Range("L75").Value = 25586.5
$ws.Range("N75").Value = -27334.5

$ws.Range("H78").Value = 25586.5
$ws.Range("J78").Value = 25586.5
$ws.Range("L78").Value = 76759.5
$ws.Range("N78").Value = -85495.5

$ws.Range("H107").Value = 38493.5
$ws.Range("J107").Value = 38493.5
$ws.Range("L107").Value = 38493.5
$ws.Range("N107").Value = -46173.5

$ws.Range("H124").Value = 20000
$ws.Range("J124").Value = 20000
$ws.Range("L124").Value = 20000
$ws.Range("N124").Value = -29820

$ws.Range("H132").Value = 19232558
$ws.Range("I132").Value = 38462756
$ws.Range("J132").Value = 2362.3076
$ws.Range("K132").Value = 115388268
$ws.Range("L132").Value = 7086.9228
$ws.Range("M132").Value = -115385738
$ws.Range("N132").Value = -12146.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 998.1429000000001
$ws.Range("I36").Value = 664.5
$ws.Range("J36").Value = 3000
$ws.Range("K36").Value = 664.5
$ws.Range("L36").Value = 3000
$ws.Range("M36").Value = -130.5
$ws.Range("N36").Value = -4068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 44313.168
$ws.Range("J41").Value = 44313.168
$ws.Range("L41").Value = 44313.168
$ws.Range("N41").Value = -45169.168

$ws.Range("H104").Value = 29997.223
$ws.Range("J104").Value = 29997.223
$ws.Range("L104").Value = 29997.223
$ws.Range("N104").Value = -35239.223

$ws.Range("H109").Value = 27116.908
$ws.Range("J109").Value = 27116.908
$ws.Range("L109").Value = 27116.908
$ws.Range("N109").Value = -29196.908

$ws.Range("H115").Value = 28925
$ws.Range("J115").Value = 28925
$ws.Range("L115").Value = 28925
$ws.Range("N115").Value = -31275

$ws.Range("H120").Value = 32610.273
$ws.Range("J120").Value = 32610.273
$ws.Range("L120").Value = 32610.273
$ws.Range("N120").Value = -39868.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2895.611
$ws.Range("I64").Value = 1655.75
$ws.Range("J64").Value = 3249.8572
$ws.Range("K64").Value = 4967.25
$ws.Range("L64").Value = 9749.571599999999
$ws.Range("M64").Value = -4697.25
$ws.Range("N64").Value = -10289.5716

$ws.Range("H67").Value = 2895.611
$ws.Range("I67").Value = 1655.75
$ws.Range("J67").Value = 3249.8572
$ws.Range("K67").Value = 4967.25
$ws.Range("L67").Value = 9749.571599999999
$ws.Range("M67").Value = -4031.25
$ws.Range("N67").Value = -11621.5716

$ws.Range("H104").Value = 3000
$ws.Range("J104").Value = 3000
$ws.Range("L104").Value = 9000
$ws.Range("N104").Value = -14242

$ws.Range("H113").Value = 3768.0938
$ws.Range("I113").Value = 7750.143
$ws.Range("J113").Value = 670.94446
$ws.Range("K113").Value = 23250.429
$ws.Range("L113").Value = 2012.83338
$ws.Range("M113").Value = -21080.429
$ws.Range("N113").Value = -6352.83338

$ws.Range("H131").Value = 2486.6711
$ws.Range("I131").Value = 7522.357
$ws.Range("J131").Value = 1349.5807
$ws.Range("K131").Value = 22567.071
$ws.Range("L131").Value = 4048.7421
$ws.Range("M131").Value = -17527.071
$ws.Range("N131").Value = -14128.7421

$ws.Range("H139").Value = 126844.32
$ws.Range("I139").Value = 347006.78
$ws.Range("J139").Value = 3002.9375
$ws.Range("K139").Value = 1041020.34
$ws.Range("L139").Value = 9008.8125
$ws.Range("M139").Value = -1035880.34
$ws.Range("N139").Value = -19288.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 42740.4
$ws.Range("J104").Value = 42740.4
$ws.Range("L104").Value = 42740.4
$ws.Range("N104").Value = -49728.4

$ws.Range("H108").Value = 20000
$ws.Range("J108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("N108").Value = -27680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 70004
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 70004
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 70004
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -70344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 24000
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H42").Value = 24000
$ws.Range("J42").Value = 24000
$ws.Range("L42").Value = 24000
$ws.Range("N42").Value = -24756

$ws.Range("H43").Value = 24712.5
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 31950
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 31950
$ws.Range("M43").Value = -2851
$ws.Range("N43").Value = -32248

$ws.Range("H120").Value = 45416
$ws.Range("J120").Value = 45416
$ws.Range("L120").Value = 45416
$ws.Range("N120").Value = -55092
